$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the email address in E2 (was "ahmed.maher@example.com")
$ws.Range("E2").Value = "Riyasmoosa@example.com"

# Turn it into a mailto hyperlink, then restore the standard Hyperlink
# cell style (Hyperlinks.Add leaves behind a tweaked style variant).
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Riyasmoosa@example.com") | Out-Null
$ws.Range("E2").Style = "Hyperlink"

# Move the active selection from J2 to H12
$ws.Range("H12").Select() | Out-Null
